$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 2912
$ws.Range("J3").Value = 3008
$ws.Range("F4").Value = 1886
$ws.Range("J4").Value = 681
$ws.Range("J6").Value = 3654
$ws.Range("F7").Value = 24076
$ws.Range("J7").Value = 10487

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 37
$ws.Range("J3").Value = 29
$ws.Range("J4").Value = 13
$ws.Range("J7").Value = 114

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 123
$ws.Range("J6").Value = 98
$ws.Range("J7").Value = 347

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 49
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 155
$ws.Range("J7").Value = 378

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 79
$ws.Range("J4").Value = 13
$ws.Range("J7").Value = 277

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 47
$ws.Range("J7").Value = 316
$ws.Range("J8").Value = 670
$ws.Range("J12").Value = 20
$ws.Range("J18").Value = 109
$ws.Range("J19").Value = 332
$ws.Range("J20").Value = 214
$ws.Range("J23").Value = 109
$ws.Range("J24").Value = 33
$ws.Range("F28").Value = 9
$ws.Range("J29").Value = 596
$ws.Range("J31").Value = 77
$ws.Range("J33").Value = 441
$ws.Range("J37").Value = 347
$ws.Range("J42").Value = 419
$ws.Range("J48").Value = 104
$ws.Range("J50").Value = 60
$ws.Range("J52").Value = 277
$ws.Range("J54").Value = 205
$ws.Range("J55").Value = 133
$ws.Range("J57").Value = 48
$ws.Range("J63").Value = 46
$ws.Range("J64").Value = 71
$ws.Range("J65").Value = 277
$ws.Range("J67").Value = 378
$ws.Range("J69").Value = 25
$ws.Range("J77").Value = 93
$ws.Range("J78").Value = 137
$ws.Range("J79").Value = 309
$ws.Range("J81").Value = 11
$ws.Range("J83").Value = 244
$ws.Range("J84").Value = 95
$ws.Range("J85").Value = 484
$ws.Range("J86").Value = 63
$ws.Range("J89").Value = 114
$ws.Range("J90").Value = 118
$ws.Range("J91").Value = 121
$ws.Range("J95").Value = 170
$ws.Range("J99").Value = 152
$ws.Range("J100").Value = 22
$ws.Range("F101").Value = 24076
$ws.Range("J101").Value = 10487

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 72
$ws.Range("J7").Value = 244

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 66
$ws.Range("J3").Value = 50
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 115
$ws.Range("J3").Value = 137
$ws.Range("J6").Value = 150
$ws.Range("J7").Value = 441

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 93
$ws.Range("J7").Value = 205

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 203
$ws.Range("J4").Value = 37
$ws.Range("J6").Value = 157
$ws.Range("J7").Value = 596

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 92
$ws.Range("J6").Value = 130
$ws.Range("J7").Value = 332

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 104

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 115
$ws.Range("J3").Value = 185
$ws.Range("J7").Value = 484

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 90
$ws.Range("J7").Value = 419

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value = 49
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 26
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J2").Value = 8
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 35
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 35
$ws.Range("J3").Value = 54
$ws.Range("J7").Value = 121

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 88
$ws.Range("J3").Value = 113
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 309

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J6").Value = 26
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 71
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 214

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J6").Value = 124
$ws.Range("J7").Value = 277

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 203
$ws.Range("J3").Value = 214
$ws.Range("J6").Value = 198
$ws.Range("J7").Value = 670

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 30
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 39
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J4").Value = 9
$ws.Range("J6").Value = 56

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J2").Value = 30
$ws.Range("J3").Value = 31
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 106
$ws.Range("J7").Value = 316

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("J3").Value = 2
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item("Sauganash,Forest Glen")
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 11

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("F4").Value = 3
$ws.Range("F7").Value = 9
